# This edit re-shuffles the per-listing data (date, variety, quality, volume,
# prices, unit, origin, $/kg, kg-or-units) across the data rows of the sheet,
# while keeping the fixed columns (Mercado ID/Mercado/Region/Codreg/Categoria
# ID/Categoria/Clasificacion) untouched. Row 2..30 hold the records; the
# mapping below says: destRow gets the data that used to live in srcRow.
# Rows not listed (7, 8, 14) keep their own data (no-ops).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 11
    3  = 12
    4  = 30
    5  = 25
    6  = 23
    9  = 2
    10 = 26
    11 = 27
    12 = 29
    13 = 5
    15 = 9
    16 = 10
    17 = 24
    18 = 28
    19 = 15
    20 = 16
    21 = 17
    22 = 4
    23 = 19
    24 = 3
    25 = 18
    26 = 22
    27 = 20
    28 = 6
    29 = 13
    30 = 21
}

# Columns that move together with a record: D (Fecha), H..Q (Variedad .. Kg o Unidades)
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17)

# 1) Snapshot the "before" values of every row that is used as a source,
#    before any writes happen (the mapping is one big permutation cycle).
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# 2) Write the snapshotted values into their destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $rowVals[$col]
    }
}
